# Apply updated coin price/volume data (GitHub Actions refresh, Wed Sep 6 09:21:44 UTC 2023).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 13/14 swap identity (WrappedEther <-> WrappedliquidstakedEther2.0) plus new price/volume.

$ws.Range("D2").Value = '25.817.00'
$ws.Range("E2").Value = '  -0.01%  '

$ws.Range("D3").Value = '1.635.36'
$ws.Range("E3").Value = '  +0.23%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '''215.37'  # force text: avoid numeric auto-coercion
$ws.Range("E5").Value = '  -0.10%  '

$ws.Range("E6").Value = '  -0.41%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("D8").Value = '''0.258'  # force text: avoid numeric auto-coercion
$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("E9").Value = '  -0.23%  '

$ws.Range("D10").Value = '''19.87'  # force text: avoid numeric auto-coercion
$ws.Range("E10").Value = '  +1.91%  '

$ws.Range("D11").Value = '''0.0781'  # force text: avoid numeric auto-coercion
$ws.Range("E11").Value = '  +0.31%  '

$ws.Range("E12").Value = '  -0.54%  '

$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '1.861.79'
$ws.Range("E13").Value = '  +0.34%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.633.78'
$ws.Range("E14").Value = '  +0.11%  '

$ws.Range("D15").Value = '''0.558'  # force text: avoid numeric auto-coercion
$ws.Range("E15").Value = '  -0.47%  '

$ws.Range("E16").Value = '  +1.29%  '

$ws.Range("D17").Value = '''63.10'  # force text: avoid numeric auto-coercion
$ws.Range("E17").Value = '  +0.08%  '

$ws.Range("D18").Value = '25.835.74'
$ws.Range("E18").Value = '  +0.02%  '

$ws.Range("E19").Value = '  +0.01%  '

$ws.Range("D20").Value = '''193.95'  # force text: avoid numeric auto-coercion
$ws.Range("E20").Value = '  -0.39%  '

$ws.Range("E21").Value = '  +1.55%  '

$ws.Range("E22").Value = '  +0.85%  '

$ws.Range("D23").Value = '''6.15'  # force text: avoid numeric auto-coercion
$ws.Range("E23").Value = '  +2.07%  '

$ws.Range("E24").Value = '  +0.05%  '

$ws.Range("D25").Value = '''1.76'  # force text: avoid numeric auto-coercion
$ws.Range("E25").Value = '  -1.81%  '

$ws.Range("D26").Value = '''140.24'  # force text: avoid numeric auto-coercion
$ws.Range("E26").Value = '  -0.64%  '

$ws.Range("E27").Value = '  -5.22%  '

$ws.Range("E28").Value = '  +1.28%  '

$ws.Range("D29").Value = '''15.49'  # force text: avoid numeric auto-coercion
$ws.Range("E29").Value = '  +0.49%  '

$ws.Range("E30").Value = '  +0.20%  '

$ws.Range("D31").Value = '''0.0494'  # force text: avoid numeric auto-coercion
$ws.Range("E31").Value = '  +1.26%  '

$ws.Range("D32").Value = '''3.32'  # force text: avoid numeric auto-coercion
$ws.Range("E32").Value = '  +1.16%  '

$ws.Range("D33").Value = '''3.25'  # force text: avoid numeric auto-coercion
$ws.Range("E33").Value = '  +0.83%  '

$ws.Range("E34").Value = '  +1.92%  '

$ws.Range("E35").Value = '  +0.23%  '

$ws.Range("D36").Value = '''0.902'  # force text: avoid numeric auto-coercion
$ws.Range("E36").Value = '  +0.64%  '

$ws.Range("E37").Value = '  -0.04%  '

$ws.Range("D38").Value = '''0.551'  # force text: avoid numeric auto-coercion
$ws.Range("E38").Value = '  +0.33%  '

$ws.Range("D39").Value = '1.112.58'
$ws.Range("E39").Value = '  -1.44%  '

$ws.Range("E40").Value = '  +0.63%  '

$ws.Range("E41").Value = '  +0.79%  '

$ws.Range("D42").Value = '''5.58'  # force text: avoid numeric auto-coercion
$ws.Range("E42").Value = '  +0.73%  '

$ws.Range("D43").Value = '''99.60'  # force text: avoid numeric auto-coercion
$ws.Range("E43").Value = '  +2.34%  '

$ws.Range("D44").Value = '''0.796'  # force text: avoid numeric auto-coercion
$ws.Range("E44").Value = '  -0.22%  '

$ws.Range("E45").Value = '  -2.16%  '

$ws.Range("D46").Value = '''55.50'  # force text: avoid numeric auto-coercion
$ws.Range("E46").Value = '  +0.32%  '

$ws.Range("E47").Value = '  +10.19%  '

$ws.Range("D48").Value = '''0.420'  # force text: avoid numeric auto-coercion
$ws.Range("E48").Value = '  -5.32%  '

$ws.Range("E49").Value = '  -0.43%  '

$ws.Range("D50").Value = '''7.64'  # force text: avoid numeric auto-coercion
$ws.Range("E50").Value = '  -0.10%  '

$ws.Range("E51").Value = '  +0.69%  '
